# Apply crypto price/volume updates described by the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force a Text number format before assigning so that numeric-looking
    # strings (e.g. "242.04") are kept as text, matching the original
    # inlineStr cell content instead of being auto-converted to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- Simple price / volume(1h) updates -----------------------------------
# Row 2
Set-TextValue $ws.Range("D2") "40.898.81"
Set-TextValue $ws.Range("E2") "  -6.98%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.196.61"
Set-TextValue $ws.Range("E3") "  -7.56%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.02%  "

# Row 5
Set-TextValue $ws.Range("D5") "242.04"
Set-TextValue $ws.Range("E5") "  +0.46%  "

# Row 6
Set-TextValue $ws.Range("E6") "  -8.00%  "

# Row 7
Set-TextValue $ws.Range("D7") "69.34"
Set-TextValue $ws.Range("E7") "  -6.93%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +0.12%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.543"
Set-TextValue $ws.Range("E9") "  -12.09%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0949"
Set-TextValue $ws.Range("E10") "  -7.84%  "

# Row 11
Set-TextValue $ws.Range("D11") "36.45"
Set-TextValue $ws.Range("E11") "  -2.96%  "

# Row 12
Set-TextValue $ws.Range("D12") "57.41"
Set-TextValue $ws.Range("E12") "  -4.84%  "

# Row 14
Set-TextValue $ws.Range("D14") "6.61"
Set-TextValue $ws.Range("E14") "  -10.06%  "

# Row 15
Set-TextValue $ws.Range("D15") "2.523.27"
Set-TextValue $ws.Range("E15") "  -7.47%  "

# Row 16
Set-TextValue $ws.Range("D16") "14.65"
Set-TextValue $ws.Range("E16") "  -11.03%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.832"
Set-TextValue $ws.Range("E17") "  -10.48%  "

# Row 18
Set-TextValue $ws.Range("D18") "2.193.85"
Set-TextValue $ws.Range("E18") "  -7.16%  "

# Row 19
Set-TextValue $ws.Range("D19") "40.846.62"
Set-TextValue $ws.Range("E19") "  -7.06%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.0₃0941"
Set-TextValue $ws.Range("E20") "  -9.40%  "

# Row 21
Set-TextValue $ws.Range("D21") "72.45"
Set-TextValue $ws.Range("E21") "  -7.38%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.05"
Set-TextValue $ws.Range("E22") "  -8.54%  "

# Row 23
Set-TextValue $ws.Range("D23") "229.73"
Set-TextValue $ws.Range("E23") "  -9.92%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.01"
Set-TextValue $ws.Range("E24") "  +6.96%  "

# Row 25
Set-TextValue $ws.Range("E25") "  +0.02%  "

# Row 26
Set-TextValue $ws.Range("D26") "3.58"
Set-TextValue $ws.Range("E26") "  -4.84%  "

# Row 27
Set-TextValue $ws.Range("E27") "  -4.03%  "

# Row 28
Set-TextValue $ws.Range("E28") "  -5.04%  "

# Row 29
Set-TextValue $ws.Range("D29") "9.67"
Set-TextValue $ws.Range("E29") "  -8.82%  "

# Row 30
Set-TextValue $ws.Range("D30") "169.56"
Set-TextValue $ws.Range("E30") "  -3.80%  "

# Row 31
Set-TextValue $ws.Range("D31") "20.24"
Set-TextValue $ws.Range("E31") "  -9.91%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.118"
Set-TextValue $ws.Range("E32") "  -9.81%  "

# Row 33
Set-TextValue $ws.Range("E33") "  -8.39%  "

# Row 35
Set-TextValue $ws.Range("D35") "5.13"
Set-TextValue $ws.Range("E35") "  -5.51%  "

# Row 36
Set-TextValue $ws.Range("E36") "  -10.58%  "

# Row 39
Set-TextValue $ws.Range("E39") "  -6.85%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0271"
Set-TextValue $ws.Range("E40") "  -3.87%  "

# Row 41
Set-TextValue $ws.Range("E41") "  -13.31%  "

# Row 42
Set-TextValue $ws.Range("D42") "62.42"
Set-TextValue $ws.Range("E42") "  -4.45%  "

# Row 43
Set-TextValue $ws.Range("D43") "4.87"
Set-TextValue $ws.Range("E43") "  -10.38%  "

# Row 44
Set-TextValue $ws.Range("E44") "  -5.99%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.192"
Set-TextValue $ws.Range("E45") "  -6.28%  "

# Row 46
Set-TextValue $ws.Range("E46") "  +0.19%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.0981"
Set-TextValue $ws.Range("E47") "  -8.99%  "

# Row 48
Set-TextValue $ws.Range("D48") "4.51"
Set-TextValue $ws.Range("E48") "  +1.39%  "

# Row 49
Set-TextValue $ws.Range("E49") "  +6.92%  "

# Row 50
Set-TextValue $ws.Range("E50") "  -6.86%  "

# Row 51
Set-TextValue $ws.Range("E51") "  -6.78%  "

# --- Rows 37/38 swap content (RenderToken <-> InjectiveProtocol) --------
# Row 37 becomes InjectiveProtocol, row 38 becomes RenderToken, each with
# their freshly updated price/volume figures.
Set-TextValue $ws.Range("B37") "InjectiveProtocol"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D37") "24.03"
Set-TextValue $ws.Range("E37") "  +15.43%  "

Set-TextValue $ws.Range("B38") "RenderToken"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D38") "3.83"
Set-TextValue $ws.Range("E38") "  -0.29%  "
